$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 2773
$ws.Range("F4").Value = 1106
$ws.Range("F5").Value = 20201
$ws.Range("F7").Value = 2401
$ws.Range("F10").Value = 465
$ws.Range("F12").Value = 259
$ws.Range("F15").Value = 386
$ws.Range("F17").Value = 488
$ws.Range("F22").Value = 109

# Sheet "演出" (Performance)
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 23
$ws.Range("F6").Value = 301
$ws.Range("F7").Value = 137
$ws.Range("F13").Value = 94
$ws.Range("F15").Value = 110

# Sheet "本地生活" (Local Life)
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 6030
$ws.Range("F3").Value = 667
$ws.Range("F4").Value = 615
$ws.Range("F5").Value = 1253

# Sheet "全部类型" (All Types)
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 6030
$ws.Range("F3").Value = 667
$ws.Range("F4").Value = 615
$ws.Range("F7").Value = 1253
$ws.Range("F8").Value = 2773
$ws.Range("F9").Value = 1106
$ws.Range("F10").Value = 20201
$ws.Range("F12").Value = 23
$ws.Range("F15").Value = 301
$ws.Range("F16").Value = 2401
$ws.Range("F18").Value = 137
$ws.Range("F20").Value = 465
$ws.Range("F22").Value = 259
$ws.Range("F28").Value = 386
$ws.Range("F32").Value = 488
$ws.Range("F33").Value = 94
$ws.Range("F37").Value = 110
$ws.Range("F38").Value = 110
$ws.Range("F49").Value = 109
